# Update the polynomial-regression result table (resultSheet!B2:L11).
# The original workbook only had data for t = 0.0 .. 0.8 (columns B..J);
# columns K and L (t = 0.9 and t = 1.0) held #NUM! errors. The fitted
# polynomials have been recomputed/extended so that every row now has
# valid values across the whole B..L range, and several previously
# computed values shifted slightly because of the extended fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..L (t = 0.0 .. 1.0), keyed by sheet row number.
$rowData = @{
    2  = @(0.7936507936507936, 0.7985175398975912, 0.7999470240275934, 0.7968687010880957, 0.7880999397162711, 0.7723411720772755, 0.7481814851195527, 0.7141198053212, 0.6686209550437986, 0.6102593894050423, 0.5379563171798417)
    3  = @(0.6627906976744186, 0.6625592712483745, 0.6577605454001819, 0.6478621622318, 0.6323665658068801, 0.610897288747511, 0.5833421956229687, 0.5500761126945722, 0.512299989687191, 0.4725673618374788, 0.4356465931125301)
    4  = @(0.5581395348837209, 0.5193952492514508, 0.47807060976512555, 0.43276393224764764, 0.38179996237145425, 0.32320010514409264, 0.25467405943894467, 0.17366269217231584, 0.07755608595454938, -0.03546866382286415, -0.16370879988422202)
    5  = @(0.25, 0.30142108200324147, 0.3527639360354631, 0.40412602621681165, 0.4556039878270791, 0.5071112501299861, 0.5580512162527507, 0.6067936876235762, 0.6499005734600232, 0.6810684722887468, 0.6898376981702953)
    6  = @(0.851063829787234, 0.8392756871350325, 0.8298855210392078, 0.8233155257857421, 0.8201712361188181, 0.8213324531964666, 0.8280820115578699, 0.842282767630592, 0.866601550374608, 0.9047093696348399, 0.9610773304917004)
    7  = @(0.8333333333333334, 0.829903494774499, 0.8315816087477425, 0.8397982646064193, 0.8559663835764175, 0.8815150662730444, 0.9179192194199888, 0.9667009099755334, 1.0293618520063987, 1.107200167791253, 1.2010108667157728)
    8  = @(0.0, -0.09999558493425532, -0.19990121713088746, -0.29968577862942636, -0.39934315648513025, -0.4988979508479157, -0.598384589357363, -0.697794786347635, -0.7969972153746345, -0.8956217756965511, -0.9928278992206233)
    9  = @(0.7, 0.7804842244320525, 0.8619827917663819, 0.9453887806576081, 1.031760661783467, 1.1224237026840023, 1.2190347556417933, 1.3234177912581395, 1.4366816849518937, 1.5567425212427146, 1.6744069080591903)
    10 = @(0.5, 0.5877714784914168, 0.6755758722216769, 0.7638295634098229, 0.852995107631085, 0.9435610710894989, 1.0360147936549124, 1.1308013356575952, 1.2282618001805785, 1.3285379784934344, 1.4313638577160734)
    11 = @(0.5, 0.5004777314014607, 0.5069132785551892, 0.519317406125767, 0.5376301729805513, 0.5616569204523728, 0.5909523851306558, 0.6247238645063495, 0.6620723359362585, 0.7036605676547337, 0.7579368499408234)
}

foreach ($rowNum in $rowData.Keys) {
    $values = $rowData[$rowNum]
    $colIndex = 2 # column B
    foreach ($val in $values) {
        $ws.Cells.Item($rowNum, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
}

Write-Host "Updated resultSheet B2:L11 with recomputed polynomial values."
